$wb = $excel.ActiveWorkbook

# --- "Chart" sheet: append the next day's row (80) ------------------------
# The date column stores plain text like "2025-12-22", not real dates. Typing
# a date-shaped string straight into Value auto-converts it to a date serial
# (and tags the cell with a date number format), so stage the text in a
# scratch cell (quote-prefixed so it stays text), copy/paste-values it into
# place, then remove the scratch cell. Paste-values carries the text over
# without pulling in the scratch cell's quote-prefix style, keeping A80 on
# the sheet's default style like every other row.
$chart = $wb.Worksheets.Item("Chart")
$scratch = $chart.Cells.Item(200, 1)
$scratch.Value = "'2025-12-23"
$scratch.Copy()
$chart.Cells.Item(80, 1).PasteSpecial(-4163)
$scratch.Delete()

$chart.Cells.Item(80, 2).Value = 0
$chart.Cells.Item(80, 3).Value = 30

# --- "Critical issues" / "Non-critical issues" header rows ----------------
# Values are unchanged ("Issue" / "Validation" / "Items"); only the shared-
# string slot they point at shifted in the source diff because the new date
# string was inserted earlier in the table. Re-asserting the same text here
# keeps both sheets' row 1 correct regardless of shared-string ordering.
foreach ($name in @("Critical issues", "Non-critical issues")) {
    $sheet = $wb.Worksheets.Item($name)
    $sheet.Cells.Item(1, 1).Value = "Issue"
    $sheet.Cells.Item(1, 2).Value = "Validation"
    $sheet.Cells.Item(1, 3).Value = "Items"
}
